# ECONOMICS EC9601A_2020.docx -- Sept 15 class update
#
# 1. Problem-set weighting text: "best 9 out of 10 ... 9% total" ->
#    "best 2 out of 10 ... 8% total".
# 2. Class-participation weight: "10%" -> "11%", and the "_GoBack" edit
#    bookmark moves from the instructor e-mail paragraph to right after
#    this new "11".
# 3. The instructor's "mailto:vaguiar@uwo.ca" link - previously a raw
#    HYPERLINK field (fldChar begin/instrText/separate/end) - becomes a
#    real w:hyperlink run, and the stray "_GoBack" bookmark that used to
#    sit in front of that field is removed (see point 2 above).

$d = $word.ActiveDocument

# --- 1. "best 9 out of 10 at 1% each, 9% total" -> "best 2 ... 8% total"
$rng = $d.Content
$rng.Find.Execute(
    "best 9 out of 10 at 1% each, 9% total",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $rng.Start
$d.Range($base + 5, $base + 6).Text  = "2"   # "9" -> "2"   (best _ out of 10)
$d.Range($base + 29, $base + 30).Text = "8"  # "9" -> "8"   (_% total)

# --- 2a. Drop the "_GoBack" bookmark currently next to the e-mail field.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2b. Replace the " HYPERLINK ""mailto:vaguiar@uwo.ca"" " field with a
#         real hyperlink run reading "vaguiar@uwo.ca".
foreach ($fld in $d.Fields) {
    if ($fld.Code.Text -match "mailto:vaguiar@uwo\.ca") {
        $insertAt = $fld.Code.Start - 1
        $fld.Delete()
        $ins = $d.Range($insertAt, $insertAt)
        $ins.Text = "vaguiar@uwo.ca"
        $linkRange = $d.Range($insertAt, $insertAt + 14)
        $d.Hyperlinks.Add($linkRange, "mailto:vaguiar@uwo.ca") | Out-Null
        break
    }
}

# --- 3. "class participation (10%)" -> "(11%)", with "_GoBack" placed
#        right after the new "11".
$rng2 = $d.Content
$rng2.Find.Execute(
    "class participation (10%)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base2 = $rng2.Start
$d.Range($base2 + 21, $base2 + 23).Text = "11"
$bmPos = $base2 + 23
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null
